$d = $word.ActiveDocument

# wdParagraph unit used with Range.Expand()/MoveEnd() etc.
$wdParagraph = 4

# --- locate the start of the block to remove -----------------------------
# The block to delete starts right after the "...Livraria da Fisica, 2005."
# paragraph (i.e. at the beginning of the blank paragraph that follows it).
# Anchor on "Nanotecnologia" (plain ASCII, unique in the document) and
# expand the found hit to its whole paragraph so Range.End lands exactly on
# the paragraph mark that ends it.
$rStart = $d.Content
$okStart = $rStart.Find.Execute("Nanotecnologia")
if (-not $okStart) {
    throw "Could not find anchor paragraph ending in 'Nanotecnologia'"
}
$rStart.Expand($wdParagraph) | Out-Null
$startPos = $rStart.End

# --- locate the end of the block to remove --------------------------------
# The block ends at the end of the copyright/footer paragraph
# ("... Creative Commons Attribution"), mark included.
$rEnd = $d.Content
$okEnd = $rEnd.Find.Execute("Creative Commons Attribution")
if (-not $okEnd) {
    throw "Could not find the footer paragraph ending in 'Creative Commons Attribution'"
}
$rEnd.Expand($wdParagraph) | Out-Null
$endPos = $rEnd.End

# --- remove the blank paragraph + the two footer paragraphs ---------------
$toRemove = $d.Range($startPos, $endPos)
$toRemove.Delete()
